$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-28 hold the "Periodo Mora" (col E) / "Valor Mora" (col F) pairs, one
# row per period. The update re-sequences them so the table reads in
# ascending period order (2001 .. 2012) followed by 2101 at the bottom,
# instead of the previous 2101, 2012, 2011 .. 2001 ordering. Each period's
# own due amount travels with it, so this is simply a top-to-bottom reversal
# of the 13-row block.
$top = 16
$bottom = 28

while ($top -lt $bottom) {
    $eTop = $ws.Range("E$top").Value()
    $eBottom = $ws.Range("E$bottom").Value()
    $ws.Range("E$top").Value = $eBottom
    $ws.Range("E$bottom").Value = $eTop

    $fTop = $ws.Range("F$top").Value()
    $fBottom = $ws.Range("F$bottom").Value()
    $ws.Range("F$top").Value = $fBottom
    $ws.Range("F$bottom").Value = $fTop

    $top++
    $bottom--
}
